$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")
$ws.Range("D13:E13").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("A8:E13").EntireRow.Delete()
Write-Output "done"
